# Refactor synthetic array: insert a new "statut_name" column (C) that
# spells out, in French, the meaning of the existing "statut_label"
# (column B) colour code. All columns from the former C (NCTId) onward
# shift one place to the right (D..M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column before the current column C (NCTId), pushing
# NCTId..intervention_type from C:L to D:M.
$ws.Columns("C:C").Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "statut_name"

# Fill in the long-form status text for every data row, derived from the
# short colour code already present in column B (statut_label).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 2).Value2

    if ($label -eq "noir") {
        $ws.Cells.Item($r, 3).Value = "pas de résultat ni de publication"
    } elseif ($label -eq "rouge") {
        $ws.Cells.Item($r, 3).Value = "résultat et / ou publication posté"
    } elseif ($label -eq "vert") {
        $ws.Cells.Item($r, 3).Value = "résultat et / ou publication posté dans les 12 mois"
    } elseif ($label -eq "orange") {
        $ws.Cells.Item($r, 3).Value = "résultat et / ou publication posté dans les 36 mois"
    }
}
